# Weekly data refresh for the Alcachofa / Terminal Hortofrutícola Agro Chillán sheet.
# A new week's record is inserted at row 15 (new Fecha/Volumen) and the
# previously-existing rows 15..20 shift down to 16..21, so row 21 is a brand
# new row at the bottom that duplicates what used to be row 20's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 21 is new: create it first, copying the static (non-changing)
#     columns from row 20, then fix up D/J/K/L/M/O/P below. --------------
# (NB: read via .Value2 — the plain .Value getter does not resolve
#  correctly as an rvalue in this host, .Value2 does.)
$ws.Range("A21").Value = $ws.Range("A20").Value2
$ws.Range("B21").Value = $ws.Range("B20").Value2
$ws.Range("C21").Value = $ws.Range("C20").Value2
$ws.Range("D21").Value = $ws.Range("D20").Value2
$ws.Range("D21").NumberFormat = $ws.Range("D20").NumberFormat
$ws.Range("E21").Value = $ws.Range("E20").Value2
$ws.Range("F21").Value = $ws.Range("F20").Value2
$ws.Range("G21").Value = $ws.Range("G20").Value2
$ws.Range("H21").Value = $ws.Range("H20").Value2
$ws.Range("I21").Value = $ws.Range("I20").Value2
$ws.Range("J21").Value = $ws.Range("J20").Value2
$ws.Range("K21").Value = $ws.Range("K20").Value2
$ws.Range("L21").Value = $ws.Range("L20").Value2
$ws.Range("M21").Value = $ws.Range("M20").Value2
$ws.Range("N21").Value = $ws.Range("N20").Value2
$ws.Range("O21").Value = $ws.Range("O20").Value2
$ws.Range("P21").Value = $ws.Range("P20").Value2
$ws.Range("Q21").Value = $ws.Range("Q20").Value2
$ws.Range("R21").Value = $ws.Range("R20").Value2

# --- Now shift the weekly figures: each row 20..16 takes on the values that
#     used to belong to the row above it (row 15..19), and row 15 gets the
#     brand-new week's Fecha/Volumen while keeping its other figures. -------
$ws.Range("D20").Value = 44425
$ws.Range("O20").Value = "Región del Maule"

$ws.Range("D19").Value = 44468
$ws.Range("J19").Value = 60
$ws.Range("K19").Value = 12000
$ws.Range("L19").Value = 13000
$ws.Range("M19").Value = 12500
$ws.Range("O19").Value = "Provincia del Elquí"
$ws.Range("P19").Value = 312

$ws.Range("D18").Value = 44435
$ws.Range("J18").Value = 120
$ws.Range("K18").Value = 14000
$ws.Range("L18").Value = 15000
$ws.Range("M18").Value = 14500
$ws.Range("P18").Value = 362

$ws.Range("D17").Value = 44453
$ws.Range("J17").Value = 160
$ws.Range("K17").Value = 12500
$ws.Range("L17").Value = 13000
$ws.Range("M17").Value = 12750
$ws.Range("P17").Value = 319

$ws.Range("D16").Value = 44475
$ws.Range("J16").Value = 120
$ws.Range("K16").Value = 11000
$ws.Range("L16").Value = 12000
$ws.Range("M16").Value = 11500
$ws.Range("P16").Value = 288

$ws.Range("D15").Value = 44488
$ws.Range("J15").Value = 100
